$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.042525944612522
$ws.Range("D2").Value = 1.043127549228891
$ws.Range("E2").Value = 1.055437458977725
$ws.Range("F2").Value = 1.062408188199581
$ws.Range("I2").Value = 1.036224564804493
$ws.Range("J2").Value = 1.047601294495798
$ws.Range("K2").Value = 1.045902290765331
$ws.Range("L2").Value = 1.05817791869641
$ws.Range("M2").Value = 1.065129616632391
$ws.Range("N2").Value = 1.04908900849553

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.044363936152056
$ws.Range("D3").Value = 1.044535441457997
$ws.Range("E3").Value = 1.057229716300947
$ws.Range("F3").Value = 1.064372322583019
$ws.Range("I3").Value = 1.036710109334959
$ws.Range("J3").Value = 1.049081348839779
$ws.Range("K3").Value = 1.047119593856062
$ws.Range("L3").Value = 1.059781122437802
$ws.Range("M3").Value = 1.066905673626248
$ws.Range("N3").Value = 1.050571164686444

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.045549153574736
$ws.Range("D4").Value = 1.045442750988007
$ws.Range("E4").Value = 1.058384543270812
$ws.Range("F4").Value = 1.065638702044659
$ws.Range("I4").Value = 1.037020939795164
$ws.Range("J4").Value = 1.050034738674543
$ws.Range("K4").Value = 1.047903068162163
$ws.Range("L4").Value = 1.060813187309106
$ws.Range("M4").Value = 1.068049954015961
$ws.Range("N4").Value = 1.051525908444128

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.046046461771276
$ws.Range("D5").Value = 1.045823315221799
$ws.Range("E5").Value = 1.058868888230019
$ws.Range("F5").Value = 1.0661700253238
$ws.Range("I5").Value = 1.037150817596677
$ws.Range("J5").Value = 1.050434530561846
$ws.Range("K5").Value = 1.048231448801049
$ws.Range("L5").Value = 1.061245818274871
$ws.Range("M5").Value = 1.068529849385963
$ws.Range("N5").Value = 1.051926268081768

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.046129906479692
$ws.Range("D6").Value = 1.045887163120285
$ws.Range("E6").Value = 1.058950145492738
$ws.Range("F6").Value = 1.066259175272466
$ws.Range("I6").Value = 1.037172578200219
$ws.Range("J6").Value = 1.050501598461705
$ws.Range("K6").Value = 1.048286527591106
$ws.Range("L6").Value = 1.061318386343155
$ws.Range("M6").Value = 1.068610358601236
$ws.Range("N6").Value = 1.051993431225738

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.045555802363934
$ws.Range("D7").Value = 1.045447839499017
$ws.Range("E7").Value = 1.058391019574734
$ws.Range("F7").Value = 1.065645805752881
$ws.Range("I7").Value = 1.037022678344909
$ws.Range("J7").Value = 1.050040084671081
$ws.Range("K7").Value = 1.047907459878959
$ws.Range("L7").Value = 1.060818973021876
$ws.Range("M7").Value = 1.068056370927043
$ws.Range("N7").Value = 1.051531262032595

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.043147961201136
$ws.Range("D8").Value = 1.043604126254273
$ws.Range("E8").Value = 1.056044186863311
$ws.Range("F8").Value = 1.063072933716577
$ws.Range("I8").Value = 1.036389354295862
$ws.Range("J8").Value = 1.048102387474065
$ws.Range("K8").Value = 1.046314562649052
$ws.Range("L8").Value = 1.058720843574798
$ws.Range("M8").Value = 1.065730882459461
$ws.Range("N8").Value = 1.049590813083302

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.038872794073152
$ws.Range("D9").Value = 1.040326324785267
$ws.Range("E9").Value = 1.051870247374375
$ws.Range("F9").Value = 1.058503212356483
$ws.Range("I9").Value = 1.035247409289139
$ws.Range("J9").Value = 1.044654170575951
$ws.Range("K9").Value = 1.043474853824204
$ws.Range("L9").Value = 1.054981908886811
$ws.Range("M9").Value = 1.061594087867622
$ws.Range("N9").Value = 1.046137699321677

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.035999666496503
$ws.Range("D10").Value = 1.038120710853524
$ws.Range("E10").Value = 1.049060147747418
$ws.Range("F10").Value = 1.055430888433171
$ws.Range("I10").Value = 1.034468249485784
$ws.Range("J10").Value = 1.04233156844864
$ws.Range("K10").Value = 1.041558749450359
$ws.Range("L10").Value = 1.052459723856508
$ws.Range("M10").Value = 1.058808469605546
$ws.Range("N10").Value = 1.043811798832939

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.034749810911255
$ws.Range("D11").Value = 1.037160599330654
$ws.Range("E11").Value = 1.047836474450026
$ws.Range("F11").Value = 1.054094040086571
$ws.Range("I11").Value = 1.034126533423219
$ws.Range("J11").Value = 1.041319965483853
$ws.Range("K11").Value = 1.040723405676879
$ws.Range("L11").Value = 1.051360245812304
$ws.Range("M11").Value = 1.057595337098881
$ws.Range("N11").Value = 1.042798759275909

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.034284665576267
$ws.Range("D12").Value = 1.03680319134745
$ws.Range("E12").Value = 1.04738088326129
$ws.Range("F12").Value = 1.053596464840359
$ws.Range("I12").Value = 1.033998945362516
$ws.Range("J12").Value = 1.040943303017719
$ws.Range("K12").Value = 1.040412254137983
$ws.Range("L12").Value = 1.050950716869762
$ws.Range("M12").Value = 1.057143652833738
$ws.Range("N12").Value = 1.04242156190587

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.034384481677231
$ws.Range("D13").Value = 1.036879892164473
$ws.Range("E13").Value = 1.047478657781366
$ws.Range("F13").Value = 1.053703242653001
$ws.Range("I13").Value = 1.034026343412425
$ws.Range("J13").Value = 1.041024139837319
$ws.Range("K13").Value = 1.040479036749065
$ws.Range("L13").Value = 1.051038613993602
$ws.Range("M13").Value = 1.057240589675839
$ws.Range("N13").Value = 1.042502513523025

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.034711380228577
$ws.Range("D14").Value = 1.03713107189242
$ws.Range("E14").Value = 1.04779883703292
$ws.Range("F14").Value = 1.054052931156597
$ws.Range("I14").Value = 1.034116000456189
$ws.Range("J14").Value = 1.041288849097338
$ws.Range("K14").Value = 1.040697703596763
$ws.Range("L14").Value = 1.051326417320296
$ws.Range("M14").Value = 1.057558022797141
$ws.Range("N14").Value = 1.042767598700557

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.034912674054162
$ws.Range("D15").Value = 1.037285728115654
$ws.Range("E15").Value = 1.047995968066308
$ws.Range("F15").Value = 1.054268250954221
$ws.Range("I15").Value = 1.03417115348034
$ws.Range("J15").Value = 1.041451824228961
$ws.Range("K15").Value = 1.040832315951941
$ws.Range("L15").Value = 1.051503591334774
$ws.Range("M15").Value = 1.057753460754267
$ws.Range("N15").Value = 1.04293080527556

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.036082491460411
$ws.Range("D16").Value = 1.038184321967493
$ws.Range("E16").Value = 1.049141211328293
$ws.Range("F16").Value = 1.055519470744803
$ws.Range("I16").Value = 1.034490836026149
$ws.Range("J16").Value = 1.04239857904956
$ws.Range("K16").Value = 1.041614067761143
$ws.Range("L16").Value = 1.052532535177652
$ws.Range("M16").Value = 1.058888832419088
$ws.Range("N16").Value = 1.043878904596599

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.036814723291158
$ws.Range("D17").Value = 1.038746616704037
$ws.Range("E17").Value = 1.049857728333324
$ws.Range("F17").Value = 1.056302562077928
$ws.Range("I17").Value = 1.03469019812643
$ws.Range("J17").Value = 1.042990858393264
$ws.Range("K17").Value = 1.042102912871944
$ws.Range("L17").Value = 1.053175975051261
$ws.Range("M17").Value = 1.059599141635935
$ws.Range("N17").Value = 1.044472025044906

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.03724126656
$ws.Range("D18").Value = 1.039074106179298
$ws.Range("E18").Value = 1.050274998689245
$ws.Range("F18").Value = 1.056758699886982
$ws.Range("I18").Value = 1.034806064993301
$ws.Range("J18").Value = 1.04333575703267
$ws.Range("K18").Value = 1.042387503250818
$ws.Range("L18").Value = 1.053550574930407
$ws.Range("M18").Value = 1.060012785284146
$ws.Range("N18").Value = 1.044817413479941

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.037386613408815
$ws.Range("D19").Value = 1.039185689494357
$ws.Range("E19").Value = 1.05041716574752
$ws.Range("F19").Value = 1.056914125851093
$ws.Range("I19").Value = 1.034845502003679
$ws.Range("J19").Value = 1.043453262879873
$ws.Range("K19").Value = 1.04248444935342
$ws.Range("L19").Value = 1.053678184633591
$ws.Range("M19").Value = 1.060153714725548
$ws.Range("N19").Value = 1.044935086198925

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.036736219357073
$ws.Range("D20").Value = 1.03868633836247
$ws.Range("E20").Value = 1.049780921522126
$ws.Range("F20").Value = 1.056218608780924
$ws.Range("I20").Value = 1.034668851721146
$ws.Range("J20").Value = 1.042927371318608
$ws.Range("K20").Value = 1.042050520886072
$ws.Range("L20").Value = 1.053107013451609
$ws.Range("M20").Value = 1.059523001458359
$ws.Range("N20").Value = 1.044408447811322

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.03461514162442
$ws.Range("D21").Value = 1.037057127401261
$ws.Range("E21").Value = 1.047704581855017
$ws.Range("F21").Value = 1.053949984754496
$ws.Range("I21").Value = 1.034089616963866
$ws.Range("J21").Value = 1.04121092406014
$ws.Range("K21").Value = 1.040633335755076
$ws.Range("L21").Value = 1.051241697946822
$ws.Range("M21").Value = 1.057464576513152
$ws.Range("N21").Value = 1.042689563000867

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.033276354178588
$ws.Range("D22").Value = 1.036028256567091
$ws.Range("E22").Value = 1.046392930701566
$ws.Range("F22").Value = 1.05251774965392
$ws.Range("I22").Value = 1.033721609226967
$ws.Range("J22").Value = 1.040126460739078
$ws.Range("K22").Value = 1.039737266219952
$ws.Range("L22").Value = 1.050062325422033
$ws.Range("M22").Value = 1.05616414077806
$ws.Range("N22").Value = 1.0416035596175

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.033986570711451
$ws.Range("D23").Value = 1.036574115340312
$ws.Range("E23").Value = 1.047088857337368
$ws.Range("F23").Value = 1.053277570927771
$ws.Range("I23").Value = 1.033917061961973
$ws.Range("J23").Value = 1.040701861585184
$ws.Range("K23").Value = 1.040212772284021
$ws.Range("L23").Value = 1.050688166034236
$ws.Range("M23").Value = 1.056854126286768
$ws.Range("N23").Value = 1.042179777598807

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.036771693620411
$ws.Range("D24").Value = 1.038713577056362
$ws.Range("E24").Value = 1.049815629258114
$ws.Range("F24").Value = 1.056256545595673
$ws.Range("I24").Value = 1.03467849853323
$ws.Range("J24").Value = 1.042956060148992
$ws.Range("K24").Value = 1.042074196251941
$ws.Range("L24").Value = 1.053138176415826
$ws.Range("M24").Value = 1.059557407989299
$ws.Range("N24").Value = 1.044437177383135

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.039981984768865
$ws.Range("D25").Value = 1.041177242181193
$ws.Range("E25").Value = 1.0529540364233
$ws.Range("F25").Value = 1.059689027979011
$ws.Range("I25").Value = 1.035545746895713
$ws.Range("J25").Value = 1.045549728851564
$ws.Range("K25").Value = 1.044212967661898
$ws.Range("L25").Value = 1.055953613663614
$ws.Range("M25").Value = 1.062668326086619
$ws.Range("N25").Value = 1.047034529392763
